# typedFullIRI.xlsx: the OTTR "iri"/"data" header pair shifted down one row.
# A7 used to hold the text "data" and is now the numeric marker 1; "iri" now
# lives at A8 (where the 1 used to be), and "data" moved down to A9 (where
# "iri" used to be).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 1
$ws.Range("A8").Value = "iri"
$ws.Range("A9").Value = "data"

# The active selection moved from A13 to A8.
$ws.Range("A8").Select()

# Column A narrowed from ~11.34 to ~8.51 stored (OOXML) character-width
# units. Range.ColumnWidth takes the pre-padding "characters" figure (the
# exporter re-adds the standard 5px/6-per-char padding), so back out that
# offset before assigning.
$ws.Columns("A").ColumnWidth = 7.671768707482998

# Cosmetic window chrome (sheet-tab/scrollbar split ratio) tweak.
$excel.ActiveWindow.TabRatio = 990
